$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 286 (shifts existing rows 286:356 down to 287:357)
$ws.Rows.Item(286).Insert()

# Populate the newly inserted row with the new price observation
$ws.Range("A286").Value = 4
$ws.Range("B286").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C286").Value = "Los Lagos"
$ws.Range("D286").Value = 44932
$ws.Range("E286").Value = 10
$ws.Range("F286").Value = 100112021
$ws.Range("G286").Value = "Ají"
$ws.Range("H286").Value = "Inferno"
$ws.Range("I286").Value = "Primera"
$ws.Range("J286").Value = 200
$ws.Range("K286").Value = 22000
$ws.Range("L286").Value = 24000
$ws.Range("M286").Value = 23000
$ws.Range("N286").Value = "$/caja 10 kilos"
$ws.Range("O286").Value = "Región de Arica y Parinacota"
$ws.Range("P286").Value = 2300
$ws.Range("Q286").Value = 10
$ws.Range("R286").Value = "Hortaliza"
